$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-20"

# Update the header label cell (I1) to match the new "through" date
$ws.Range("I1").Value = "2022 (through 03-20)"

# Update March's "Total" column (row 4) value
$ws.Range("I4").Value = 88

# Update the grand Total row (row 14) "Total" column value
$ws.Range("I14").Value = 388
